$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.ClearContents()
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.989.04'
Set-TextValue $ws.Range('E2') '  -0.55%  '
Set-TextValue $ws.Range('D3') '1.859.53'
Set-TextValue $ws.Range('E3') '  -1.02%  '
Set-TextValue $ws.Range('E4') '  +0.19%  '
Set-TextValue $ws.Range('D5') '312.13'
Set-TextValue $ws.Range('E5') '  -0.46%  '
Set-TextValue $ws.Range('E6') '  +0.13%  '
Set-TextValue $ws.Range('D7') '0.5137'
Set-TextValue $ws.Range('E7') '  +0.85%  '
Set-TextValue $ws.Range('D8') '0.3830'
Set-TextValue $ws.Range('E8') '  -0.53%  '
Set-TextValue $ws.Range('D9') '0.08230'
Set-TextValue $ws.Range('E9') '  -9.94%  '
Set-TextValue $ws.Range('D10') '1.109'
Set-TextValue $ws.Range('E10') '  -1.18%  '
Set-TextValue $ws.Range('E11') '  -0.31%  '
Set-TextValue $ws.Range('D12') '6.192'
Set-TextValue $ws.Range('E12') '  -2.70%  '
Set-TextValue $ws.Range('D13') '20.55'
Set-TextValue $ws.Range('E13') '  -1.12%  '
Set-TextValue $ws.Range('D14') '1.863.35'
Set-TextValue $ws.Range('E14') '  -0.38%  '
Set-TextValue $ws.Range('D15') '7.263'
Set-TextValue $ws.Range('E15') '  +0.79%  '
Set-TextValue $ws.Range('E16') '  +0.06%  '
Set-TextValue $ws.Range('D17') '0.00001095'
Set-TextValue $ws.Range('E17') '  -1.90%  '
Set-TextValue $ws.Range('D18') '90.59'
Set-TextValue $ws.Range('E18') '  -0.70%  '
Set-TextValue $ws.Range('D19') '0.06646'
Set-TextValue $ws.Range('E19') '  +0.41%  '
Set-TextValue $ws.Range('D20') '17.64'
Set-TextValue $ws.Range('E20') '  -3.17%  '
Set-TextValue $ws.Range('D22') '6.006'
Set-TextValue $ws.Range('E22') '  -1.85%  '
Set-TextValue $ws.Range('D23') '28.021.37'
Set-TextValue $ws.Range('E23') '  -0.55%  '
Set-TextValue $ws.Range('D24') '11.06'
Set-TextValue $ws.Range('E24') '  -3.42%  '
Set-TextValue $ws.Range('D25') '2.267'
Set-TextValue $ws.Range('E25') '  -0.55%  '
Set-TextValue $ws.Range('D26') '2.075.76'
Set-TextValue $ws.Range('E26') '  -0.46%  '
Set-TextValue $ws.Range('D27') '2.511'
Set-TextValue $ws.Range('E27') '  -2.55%  '
Set-TextValue $ws.Range('D28') '157.64'
Set-TextValue $ws.Range('E28') '  +0.17%  '
Set-TextValue $ws.Range('D29') '20.45'
Set-TextValue $ws.Range('E29') '  -1.72%  '
Set-TextValue $ws.Range('D30') '124.66'
Set-TextValue $ws.Range('E30') '  -1.67%  '
Set-TextValue $ws.Range('E31') '  +1.15%  '
Set-TextValue $ws.Range('E32') '  -3.19%  '
Set-TextValue $ws.Range('D33') '5.899'
Set-TextValue $ws.Range('E33') '  +4.97%  '
Set-TextValue $ws.Range('D34') '3.594'
Set-TextValue $ws.Range('D35') '9.390'
Set-TextValue $ws.Range('E35') '  -3.31%  '
Set-TextValue $ws.Range('D36') '0.02415'
Set-TextValue $ws.Range('E36') '  -3.29%  '
Set-TextValue $ws.Range('D37') '0.06499'
Set-TextValue $ws.Range('E37') '  -1.19%  '
Set-TextValue $ws.Range('E38') '  +0.09%  '
Set-TextValue $ws.Range('D39') '0.6539'
Set-TextValue $ws.Range('E39') '  +1.92%  '
Set-TextValue $ws.Range('D40') '1.196'
Set-TextValue $ws.Range('E40') '  -1.46%  '
Set-TextValue $ws.Range('D41') '4.991'
Set-TextValue $ws.Range('E41') '  +1.38%  '
Set-TextValue $ws.Range('D42') '1.215'
Set-TextValue $ws.Range('E42') '  -2.40%  '
Set-TextValue $ws.Range('D43') '11.17'
Set-TextValue $ws.Range('E43') '  -3.45%  '
Set-TextValue $ws.Range('D44') '0.6166'
Set-TextValue $ws.Range('E44') '  +2.45%  '
Set-TextValue $ws.Range('D45') '13.03'
Set-TextValue $ws.Range('E45') '  -1.18%  '
Set-TextValue $ws.Range('D46') '1.283'
Set-TextValue $ws.Range('E46') '  +0.60%  '
Set-TextValue $ws.Range('E47') '  -0.19%  '
Set-TextValue $ws.Range('D48') '2.006'
Set-TextValue $ws.Range('E48') '  +0.25%  '
Set-TextValue $ws.Range('D49') '1.214'
Set-TextValue $ws.Range('E49') '  -1.52%  '
Set-TextValue $ws.Range('D50') '120.82'
Set-TextValue $ws.Range('E50') '  -0.56%  '
Set-TextValue $ws.Range('D51') '78.18'
Set-TextValue $ws.Range('E51') '  -2.18%  '
